$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 44841
$ws.Range("I6").Value = 50414.875
$ws.Range("J6").Value = 250
$ws.Range("K6").Value = 151244.625
$ws.Range("L6").Value = 750
$ws.Range("M6").Value = -151132.625
$ws.Range("N6").Value = -974

# Row 38
$ws.Range("H38").Value = 6976.4546
$ws.Range("I38").Value = 217.625
$ws.Range("K38").Value = 652.875
$ws.Range("M38").Value = -280.875

# Row 100
$ws.Range("H100").Value = 3538.5
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3538.5
$ws.Range("K100").Value = 0
$ws.Range("M100").Value = 3538.5
$ws.Range("N100").Value = -4620.5
$ws.Range("L100").ClearContents()

# Row 135
$ws.Range("H135").Value = 482.57895
$ws.Range("I135").Value = 471.70587
$ws.Range("J135").Value = 575
$ws.Range("K135").Value = 4245.35283
$ws.Range("L135").Value = 5175
$ws.Range("M135").Value = -1710.35283
$ws.Range("N135").Value = -10245

# Row 138
$ws.Range("H138").Value = 2304.275
$ws.Range("I138").Value = 2265.8235
$ws.Range("J138").Value = 2332.6956
$ws.Range("K138").Value = 6797.470499999999
$ws.Range("L138").Value = 6998.0868
$ws.Range("M138").Value = -1657.470499999999
$ws.Range("N138").Value = -17278.0868

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4119.0977
$ws.Range("I61").Value = 3683.6765
$ws.Range("K61").Value = 3683.6765
$ws.Range("M61").Value = -3471.6765

# Row 136
$ws.Range("H136").Value = 4119.0977
$ws.Range("I136").Value = 3683.6765
$ws.Range("K136").Value = 11051.0295
$ws.Range("M136").Value = -8501.029500000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2936.077
$ws.Range("I20").Value = 1296.5
$ws.Range("J20").Value = 5559.4
$ws.Range("K20").Value = 1296.5
$ws.Range("L20").Value = 5559.4
$ws.Range("M20").Value = -1049.5
$ws.Range("N20").Value = -6053.4

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 11639.75
$ws.Range("I2").Value = 3279.5
$ws.Range("J2").Value = 20000
$ws.Range("K2").Value = 3279.5
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = -3166.5
$ws.Range("N2").Value = -20226

# Row 7
$ws.Range("H7").Value = 856.2
$ws.Range("I7").Value = 854.8125
$ws.Range("K7").Value = 854.8125
$ws.Range("M7").Value = -741.8125

# Row 16
$ws.Range("H16").Value = 1056
$ws.Range("I16").Value = 907.1818
$ws.Range("K16").Value = 907.1818
$ws.Range("M16").Value = -620.1818

# Row 41
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("N41").Value = 0
$ws.Range("L41").ClearContents()

# Row 51
$ws.Range("H51").Value = 14000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

# Row 59
$ws.Range("H59").Value = 124999.5
$ws.Range("J59").Value = 124999.5
$ws.Range("L59").Value = 124999.5
$ws.Range("N59").Value = -127289.5

# Row 60
$ws.Range("H60").Value = 38517.5
$ws.Range("J60").Value = 40384.25
$ws.Range("L60").Value = 40384.25
$ws.Range("N60").Value = -41406.25

# Row 61
$ws.Range("H61").Value = 14000
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# Row 62
$ws.Range("H62").Value = 6852.7144
$ws.Range("I62").Value = 5761.7144
$ws.Range("J62").Value = 7398.2144
$ws.Range("K62").Value = 5761.7144
$ws.Range("L62").Value = 7398.2144
$ws.Range("M62").Value = -5137.7144
$ws.Range("N62").Value = -8646.214400000001

# Row 65
$ws.Range("H65").Value = 6852.7144
$ws.Range("I65").Value = 5761.7144
$ws.Range("J65").Value = 7398.2144
$ws.Range("K65").Value = 28808.572
$ws.Range("L65").Value = 36991.072
$ws.Range("M65").Value = -25688.572
$ws.Range("N65").Value = -43231.072

# Row 107
$ws.Range("H107").Value = 4005.8462
$ws.Range("I107").Value = 3839.6667
$ws.Range("J107").Value = 6000
$ws.Range("K107").Value = 3839.6667
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = -1919.6667
$ws.Range("N107").Value = -9840

# Row 113
$ws.Range("H113").Value = 1056
$ws.Range("I113").Value = 907.1818
$ws.Range("K113").Value = 907.1818
$ws.Range("M113").Value = 1262.8182

# Row 122
$ws.Range("H122").Value = 4251.7334
$ws.Range("I122").Value = 3195.889
$ws.Range("K122").Value = 9587.667000000001
$ws.Range("M122").Value = -7137.667000000001

# Row 134
$ws.Range("H134").Value = 2103.1936
$ws.Range("I134").Value = 1479.619
$ws.Range("J134").Value = 3412.7
$ws.Range("K134").Value = 4438.857
$ws.Range("L134").Value = 10238.1
$ws.Range("M134").Value = -1903.857
$ws.Range("N134").Value = -15308.1

$ws = $wb.Worksheets.Item("CUL")
# Row 19
$ws.Range("H19").Value = 400
$ws.Range("I19").Value = 400
$ws.Range("K19").Value = 1200
$ws.Range("M19").Value = -1026

# Row 22
$ws.Range("H22").Value = 1482.1
$ws.Range("I22").Value = 3000
$ws.Range("K22").Value = 9000
$ws.Range("M22").Value = -8831

# Row 25
$ws.Range("H25").Value = 422.5
$ws.Range("I25").Value = 422.5
$ws.Range("K25").Value = 1267.5
$ws.Range("M25").Value = -1098.5

# Row 27
$ws.Range("H27").Value = 1482.1
$ws.Range("I27").Value = 3000
$ws.Range("K27").Value = 9000
$ws.Range("M27").Value = -8898

# Row 30
$ws.Range("H30").Value = 422.5
$ws.Range("I30").Value = 422.5
$ws.Range("K30").Value = 1267.5
$ws.Range("M30").Value = -1165.5

# Row 44
$ws.Range("H44").Value = 834.8333
$ws.Range("I44").Value = 913
$ws.Range("J44").Value = 444
$ws.Range("K44").Value = 2739
$ws.Range("L44").Value = 1332
$ws.Range("M44").Value = -2341
$ws.Range("N44").Value = -2128

# Row 47
$ws.Range("H47").Value = 443.25
$ws.Range("I47").Value = 443
$ws.Range("K47").Value = 1329
$ws.Range("M47").Value = -898

# Row 131
$ws.Range("H131").Value = 3199.5857
$ws.Range("I131").Value = 1658.2222
$ws.Range("J131").Value = 4167.4185
$ws.Range("K131").Value = 4974.6666
$ws.Range("L131").Value = 12502.2555
$ws.Range("M131").Value = 65.33340000000044
$ws.Range("N131").Value = -22582.2555

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 39994.5
$ws.Range("J46").Value = 39994.5
$ws.Range("L46").Value = 39994.5
$ws.Range("N46").Value = -40306.5

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1187.3914
$ws.Range("I16").Value = 1147.7894
$ws.Range("K16").Value = 1147.7894
$ws.Range("M16").Value = -977.7893999999999

# Row 22
$ws.Range("H22").Value = 1370.7
$ws.Range("I22").Value = 1401.2778
$ws.Range("J22").Value = 1324.8334
$ws.Range("K22").Value = 1401.2778
$ws.Range("L22").Value = 1324.8334
$ws.Range("M22").Value = -1106.2778
$ws.Range("N22").Value = -1914.8334

# Row 27
$ws.Range("H27").Value = 1370.7
$ws.Range("I27").Value = 1401.2778
$ws.Range("J27").Value = 1324.8334
$ws.Range("K27").Value = 1401.2778
$ws.Range("L27").Value = 1324.8334
$ws.Range("M27").Value = -1294.2778
$ws.Range("N27").Value = -1538.8334

# Row 46
$ws.Range("H46").Value = 2937.2104
$ws.Range("I46").Value = 936
$ws.Range("K46").Value = 936
$ws.Range("M46").Value = -748

# Row 68
$ws.Range("H68").Value = 5954.5454
$ws.Range("I68").Value = 4078.3076
$ws.Range("J68").Value = 8664.666999999999
$ws.Range("K68").Value = 4078.3076
$ws.Range("L68").Value = 8664.666999999999
$ws.Range("M68").Value = -3329.3076
$ws.Range("N68").Value = -10162.667

# Row 71
$ws.Range("H71").Value = 5954.5454
$ws.Range("I71").Value = 4078.3076
$ws.Range("J71").Value = 8664.666999999999
$ws.Range("K71").Value = 20391.538
$ws.Range("L71").Value = 43323.335
$ws.Range("M71").Value = -16647.538
$ws.Range("N71").Value = -50811.335

# Row 108
$ws.Range("H108").Value = 62623
$ws.Range("J108").Value = 62623
$ws.Range("L108").Value = 62623
$ws.Range("N108").Value = -70303

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 3588.476
$ws.Range("I126").Value = 3642.5881
$ws.Range("K126").Value = 10927.7643
$ws.Range("M126").Value = -8457.764299999999

# Row 136
$ws.Range("H136").Value = 4180.082
$ws.Range("I136").Value = 2687.4048
$ws.Range("K136").Value = 8062.214399999999
$ws.Range("M136").Value = -5512.214399999999

